$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.649.27"
$ws.Range("E2").Value = "  -1.53%  "
$ws.Range("D3").Value = "1.669.08"
$ws.Range("E3").Value = "  -3.17%  "
$ws.Range("E4").Value = "  -0.05%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "215.27"
$c.ClearFormats()
$ws.Range("E5").Value = "  -1.65%  "
$ws.Range("E6").Value = "  -2.34%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -1.36%  "
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("E10").Value = "  -1.80%  "
$ws.Range("E11").Value = "  -2.03%  "
$ws.Range("D12").Value = "1.904.39"
$ws.Range("E12").Value = "  -3.19%  "
$ws.Range("D13").Value = "1.666.99"
$ws.Range("E13").Value = "  -3.39%  "
$ws.Range("E14").Value = "  -3.10%  "
$ws.Range("E15").Value = "  -0.53%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "66.77"
$c.ClearFormats()
$ws.Range("D17").Value = "27.634.91"
$ws.Range("E17").Value = "  -1.43%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "243.86"
$c.ClearFormats()
$ws.Range("E18").Value = "  +0.40%  "
$ws.Range("D19").Value = "0.0₃0732"
$ws.Range("E19").Value = "  -3.20%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.68"
$c.ClearFormats()
$ws.Range("E20").Value = "  -4.62%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("E22").Value = "  -2.85%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "9.35"
$c.ClearFormats()
$ws.Range("E23").Value = "  -3.84%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.05"
$c.ClearFormats()
$ws.Range("E24").Value = "  -3.76%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "147.43"
$c.ClearFormats()
$ws.Range("E25").Value = "  -1.09%  "
$ws.Range("E26").Value = "  -3.69%  "
$ws.Range("E27").Value = "  -1.17%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("E29").Value = "  -2.30%  "
$ws.Range("E30").Value = "  +3.04%  "
$ws.Range("E31").Value = "  -1.28%  "
$ws.Range("E32").Value = "  -2.37%  "
$ws.Range("D33").Value = "1.471.35"
$ws.Range("E33").Value = "  -1.72%  "
$ws.Range("E34").Value = "  -4.77%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.58"
$c.ClearFormats()
$ws.Range("E35").Value = "  -4.97%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.931"
$c.ClearFormats()
$ws.Range("E36").Value = "  -2.77%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "2.37"
$c.ClearFormats()
$ws.Range("E37").Value = "  -1.11%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.577"
$c.ClearFormats()
$ws.Range("E38").Value = "  -4.93%  "
$ws.Range("E39").Value = "  -1.58%  "
$ws.Range("E40").Value = "  -1.69%  "
$ws.Range("E41").Value = "  -4.57%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.ClearFormats()
$ws.Range("E42").Value = "  -0.06%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "5.43"
$c.ClearFormats()
$ws.Range("E43").Value = "  -7.19%  "
$ws.Range("E44").Value = "  -2.89%  "
$ws.Range("D45").Value = "1.812.53"
$ws.Range("E45").Value = "  -3.13%  "
$ws.Range("E46").Value = "  -1.87%  "
$ws.Range("E47").Value = "  -0.83%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "89.34"
$c.ClearFormats()
$ws.Range("E48").Value = "  -1.73%  "
$ws.Range("E49").Value = "  -4.33%  "
$ws.Range("E50").Value = "  -2.05%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "7.89"
$c.ClearFormats()
